# Edit script: rebuild the data table with a new set of accounts,
# dropping the "password" and "remember_token" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the old columns G (password) and H (remember_token) ---
$ws.Range("G1:H6").Clear()

# --- Clear the old data area so we can rewrite it cleanly ---
$ws.Range("A1:H6").ClearContents()

# --- Header row ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "email_verified_at"
$ws.Range("E1").Value = "created_at"
$ws.Range("F1").Value = "updated_at"

# --- Data rows ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Maze Clarion Badlon"
$ws.Range("C2").Value = "badlonmazeclarion@gmail.com"
$ws.Range("E2").Value = "2023-06-09T09:20:30.000000Z"
$ws.Range("F2").Value = "2023-06-09T09:20:30.000000Z"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Maze Clarion Badlon"
$ws.Range("C3").Value = "badlonmaze@gmail.com"
$ws.Range("E3").Value = "2023-06-09T11:08:16.000000Z"
$ws.Range("F3").Value = "2023-06-09T11:08:16.000000Z"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Maze Clarion Badlon"
$ws.Range("C4").Value = "badlon@gmail.com"
$ws.Range("E4").Value = "2023-06-09T11:09:40.000000Z"
$ws.Range("F4").Value = "2023-06-09T11:09:40.000000Z"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Mazeu"
$ws.Range("C5").Value = "maze@gmail.com"
$ws.Range("E5").Value = "2023-06-09T11:25:51.000000Z"
$ws.Range("F5").Value = "2023-06-09T11:25:51.000000Z"

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "Kirk J-Son Matic"
$ws.Range("C6").Value = "kirk@gmail.com"
$ws.Range("E6").Value = "2023-06-09T10:35:08.000000Z"
$ws.Range("F6").Value = "2023-06-09T10:35:08.000000Z"

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Rodel Cuyag"
$ws.Range("C7").Value = "rodel@gmail.com"
$ws.Range("E7").Value = "2023-06-11T01:05:32.000000Z"
$ws.Range("F7").Value = "2023-06-11T01:05:32.000000Z"
